$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 543.1429000000001
$ws.Range("I19").Value = 288.7
$ws.Range("J19").Value = 644.92
$ws.Range("K19").Value = 288.7
$ws.Range("L19").Value = 644.92
$ws.Range("M19").Value = -113.7
$ws.Range("N19").Value = -994.92

$ws.Range("H32").Value = 983.4
$ws.Range("I32").Value = 639
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 639
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = -313
$ws.Range("N32").Value = -2152

$ws.Range("H98").Value = 810.3200000000001
$ws.Range("I98").Value = 760.5
$ws.Range("J98").Value = 2006
$ws.Range("K98").Value = 760.5
$ws.Range("L98").Value = 2006
$ws.Range("M98").Value = 737.5
$ws.Range("N98").Value = -5002

$ws.Range("H107").Value = 905.8461
$ws.Range("I107").Value = 807
$ws.Range("J107").Value = 1235.3334
$ws.Range("K107").Value = 807
$ws.Range("L107").Value = 1235.3334
$ws.Range("M107").Value = 1113
$ws.Range("N107").Value = -5075.3334

$ws.Range("H122").Value = 810.3200000000001
$ws.Range("I122").Value = 760.5
$ws.Range("J122").Value = 2006
$ws.Range("K122").Value = 2281.5
$ws.Range("L122").Value = 6018
$ws.Range("M122").Value = 168.5
$ws.Range("N122").Value = -10918

$ws.Range("H138").Value = 3111.9148
$ws.Range("I138").Value = 1747.258
$ws.Range("K138").Value = 5241.774
$ws.Range("M138").Value = -101.7740000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 954.3333
$ws.Range("I2").Value = 970.129
$ws.Range("J2").Value = 909.8182
$ws.Range("K2").Value = 970.129
$ws.Range("L2").Value = 909.8182
$ws.Range("M2").Value = -857.129
$ws.Range("N2").Value = -1135.8182

$ws.Range("H32").Value = 7339
$ws.Range("I32").Value = 5890.521
$ws.Range("J32").Value = 17271.428
$ws.Range("K32").Value = 5890.521
$ws.Range("L32").Value = 17271.428
$ws.Range("M32").Value = -5603.521
$ws.Range("N32").Value = -17845.428

$ws.Range("H35").Value = 10094.889
$ws.Range("I35").Value = 2170.8
$ws.Range("K35").Value = 2170.8
$ws.Range("M35").Value = -1764.8

$ws.Range("H74").Value = 4991.6587
$ws.Range("I74").Value = 2513.0571
$ws.Range("K74").Value = 2513.0571
$ws.Range("M74").Value = -1639.0571

$ws.Range("H77").Value = 4991.6587
$ws.Range("I77").Value = 2513.0571
$ws.Range("K77").Value = 12565.2855
$ws.Range("M77").Value = -8197.2855

$ws.Range("H116").Value = 954.3333
$ws.Range("I116").Value = 970.129
$ws.Range("J116").Value = 909.8182
$ws.Range("K116").Value = 970.129
$ws.Range("L116").Value = 909.8182
$ws.Range("M116").Value = 1323.871
$ws.Range("N116").Value = -5497.8182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 954.3333
$ws.Range("I3").Value = 970.129
$ws.Range("J3").Value = 909.8182
$ws.Range("K3").Value = 970.129
$ws.Range("L3").Value = 909.8182
$ws.Range("M3").Value = -856.129
$ws.Range("N3").Value = -1137.8182

$ws.Range("H56").Value = 35777.5
$ws.Range("J56").Value = 35777.5
$ws.Range("L56").Value = 35777.5
$ws.Range("N56").Value = -37255.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 13802.444
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = $null

$ws.Range("H58").Value = 1685783.9
$ws.Range("I58").Value = 2675230.5
$ws.Range("J58").Value = 3724.6
$ws.Range("K58").Value = 2675230.5
$ws.Range("L58").Value = 3724.6
$ws.Range("M58").Value = -2675027.5
$ws.Range("N58").Value = -4130.6

$ws.Range("H136").Value = 1685783.9
$ws.Range("I136").Value = 2675230.5
$ws.Range("J136").Value = 3724.6
$ws.Range("K136").Value = 8025691.5
$ws.Range("L136").Value = 11173.8
$ws.Range("M136").Value = -8023141.5
$ws.Range("N136").Value = -16273.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 205
$ws.Range("I5").Value = 205
$ws.Range("K5").Value = 205
$ws.Range("M5").Value = -93

$ws.Range("H70").Value = 6530.4346

$ws.Range("H73").Value = 6530.4346

$ws.Range("H126").Value = 3061.1875
$ws.Range("I126").Value = 1697.9
$ws.Range("J126").Value = 5333.3335
$ws.Range("K126").Value = 5093.700000000001
$ws.Range("L126").Value = 16000.0005
$ws.Range("M126").Value = -2623.700000000001
$ws.Range("N126").Value = -20940.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2050857.2
$ws.Range("I2").Value = 7001500.5
$ws.Range("J2").Value = 70600
$ws.Range("K2").Value = 7001500.5
$ws.Range("L2").Value = 70600
$ws.Range("M2").Value = -7001388.5
$ws.Range("N2").Value = -70824

$ws.Range("H22").Value = 603.26666
$ws.Range("I22").Value = 622.1111
$ws.Range("J22").Value = 575
$ws.Range("K22").Value = 622.1111
$ws.Range("L22").Value = 575
$ws.Range("M22").Value = -327.1111
$ws.Range("N22").Value = -1165

$ws.Range("H27").Value = 603.26666
$ws.Range("I27").Value = 622.1111
$ws.Range("J27").Value = 575
$ws.Range("K27").Value = 622.1111
$ws.Range("L27").Value = 575
$ws.Range("M27").Value = -515.1111
$ws.Range("N27").Value = -789

$ws.Range("H132").Value = 4188.8213
$ws.Range("I132").Value = 3658.476
$ws.Range("J132").Value = 5779.857
$ws.Range("K132").Value = 10975.428
$ws.Range("L132").Value = 17339.571
$ws.Range("M132").Value = -8445.428
$ws.Range("N132").Value = -22399.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 14844.615
$ws.Range("I2").Value = 80000
$ws.Range("J2").Value = 2998.182
$ws.Range("K2").Value = 80000
$ws.Range("L2").Value = 2998.182
$ws.Range("M2").Value = -79888
$ws.Range("N2").Value = -3222.182

$ws.Range("H113").Value = 980.75
$ws.Range("I113").Value = 433.7647
$ws.Range("J113").Value = 1385.0435
$ws.Range("K113").Value = 1301.2941
$ws.Range("L113").Value = 4155.1305
$ws.Range("M113").Value = 868.7058999999999
$ws.Range("N113").Value = -8495.130499999999

$ws.Range("H132").Value = 6214.7144
$ws.Range("I132").Value = 9334.666999999999
$ws.Range("J132").Value = 3874.75
$ws.Range("K132").Value = 28004.001
$ws.Range("L132").Value = 11624.25
$ws.Range("M132").Value = -25474.001
$ws.Range("N132").Value = -16684.25

$ws.Range("H136").Value = 3807.5781
$ws.Range("I136").Value = 1626
$ws.Range("J136").Value = 6440.517
$ws.Range("K136").Value = 4878
$ws.Range("L136").Value = 19321.551
$ws.Range("M136").Value = -2328
$ws.Range("N136").Value = -24421.551
